# Update "想去人数" (want-to-go count) figures in column F for several
# conghua/comic-convention events. The same events appear on the
# "展览" sheet and on the "全部类型" sheet (which has one extra leading
# row), so each sheet is updated independently with its own row numbers.

$wb = $excel.ActiveWorkbook

# -- Sheet "展览" --------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 22
$ws1.Range("F4").Value  = 1375
$ws1.Range("F6").Value  = 1031
$ws1.Range("F7").Value  = 10678
$ws1.Range("F10").Value = 291
$ws1.Range("F11").Value = 1034
$ws1.Range("F12").Value = 703
$ws1.Range("F13").Value = 12042
$ws1.Range("F16").Value = 121

# -- Sheet "全部类型" ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 22
$ws4.Range("F5").Value  = 1375
$ws4.Range("F7").Value  = 1031
$ws4.Range("F8").Value  = 10678
$ws4.Range("F11").Value = 291
$ws4.Range("F12").Value = 1034
$ws4.Range("F13").Value = 703
$ws4.Range("F14").Value = 12042
$ws4.Range("F17").Value = 121
